# Issues #2, #3, #51, #54, #55 Basic config set and send
#
# Mark the listed issues "DONE" (Status column, D) and log a new issue
# (#66) "handle error in Repository services", also already DONE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Mark existing issues as DONE (column D = Status)
$ws.Range("D3").Value = "DONE"
$ws.Range("D4").Value = "DONE"
$ws.Range("D51").Value = "DONE"
$ws.Range("D54").Value = "DONE"
$ws.Range("D55").Value = "DONE"
$ws.Range("D57").Value = "DONE"

# Log new issue #66 and mark it DONE too
$ws.Range("A66").Value = 66
$ws.Range("D66").Value = "DONE"
$ws.Range("F66").Value = "handle error in Repository services"
$ws.Rows(66).RowHeight = 29

# Leave the view where the user last left it, on the new status cell
$ws.Range("D57").Select() | Out-Null
